$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.818.77'
$ws.Range("E2").Value = '  -0.65%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.631.45'
$ws.Range("E3").Value = '  -0.68%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.63'
$ws.Range("E5").Value = '  +0.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5062'
$ws.Range("E6").Value = '  -0.11%  '

$ws.Range("E7").Value = '  +0.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2577'
$ws.Range("E8").Value = '  -0.21%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06433'
$ws.Range("E9").Value = '  +0.98%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.48'
$ws.Range("E10").Value = '  -2.10%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07788'
$ws.Range("E11").Value = '  +0.61%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.271'
$ws.Range("E12").Value = '  -0.74%  '

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.854.79'
$ws.Range("E13").Value = '  -0.77%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.627.67'
$ws.Range("E14").Value = '  -0.74%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5604'
$ws.Range("E15").Value = '  +2.19%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅7577'
$ws.Range("E16").Value = '  -2.31%  '

$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.97'
$ws.Range("E17").Value = '  -2.18%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.824.78'
$ws.Range("E18").Value = '  -0.77%  '

$ws.Range("E19").Value = '  -0.09%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '194.31'
$ws.Range("E20").Value = '  -1.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.328'
$ws.Range("E21").Value = '  -3.27%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.845'
$ws.Range("E22").Value = '  -1.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.016'
$ws.Range("E23").Value = '  -2.37%  '

$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.797'
$ws.Range("E25").Value = '  -5.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.11'
$ws.Range("E26").Value = '  -1.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1279'
$ws.Range("E27").Value = '  +1.20%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.748'
$ws.Range("E28").Value = '  -1.93%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.41'
$ws.Range("E29").Value = '  -1.67%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.240'
$ws.Range("E30").Value = '  -0.11%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04870'
$ws.Range("E31").Value = '  -0.64%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.286'
$ws.Range("E32").Value = '  +0.20%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.216'
$ws.Range("E33").Value = '  +0.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.551'
$ws.Range("E34").Value = '  -0.45%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.381'
$ws.Range("E35").Value = '  +0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8962'
$ws.Range("E36").Value = '  -2.63%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.567'
$ws.Range("E37").Value = '  -0.22%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.128.35'
$ws.Range("E38").Value = '  +0.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5490'
$ws.Range("E39").Value = '  -1.21%  '

$ws.Range("E40").Value = '  -0.63%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9937'
$ws.Range("E41").Value = '  -0.84%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.537'
$ws.Range("E42").Value = '  -1.44%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7978'
$ws.Range("E43").Value = '  -0.92%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.28'
$ws.Range("E44").Value = '  -1.40%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.780.91'
$ws.Range("E45").Value = '  +0.06%  '

$ws.Range("E46").Value = '  -4.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4444'
$ws.Range("E47").Value = '  -1.88%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.31'
$ws.Range("E48").Value = '  -0.15%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05057'
$ws.Range("E49").Value = '  -2.56%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.664'
$ws.Range("E50").Value = '  +0.68%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9974'
$ws.Range("E51").Value = '  -0.49%  '
